# Harmonogram_aktivita_MD.xlsx - apply tracked-change edit
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "List1"
$ws2 = $wb.Worksheets.Item(2)   # "Mereni aktivity"

# --- List1: updated measurement timestamps (dependent formulas recalc automatically) ---
$ws1.Range("A37").Value = 45516.675694444442
$ws1.Range("A40").Value = 45519.638888888891
$ws1.Range("A44").Value = 45523.673611111109

# --- Mereni aktivity: fill previously-empty rows 20-26 of Tabulka1 with measurement data ---
# Row 20
$ws2.Cells.Item(20, 1).Value = 45502.631944444445
$ws2.Cells.Item(20, 2).Value = 0.013
$ws2.Cells.Item(20, 3).Value = 0.011
$ws2.Cells.Item(20, 4).Value = 0.011
$ws2.Cells.Item(20, 5).Value = 0.011
$ws2.Cells.Item(20, 6).Value = 0.012
$ws2.Cells.Item(20, 7).Value = 115.2
$ws2.Cells.Item(20, 8).Value = 115.1
$ws2.Cells.Item(20, 9).Value = 115.1
$ws2.Cells.Item(20, 10).Value = 115.1
$ws2.Cells.Item(20, 11).Value = 115.0
$ws2.Cells.Item(20, 12).Value = 115.0
$ws2.Cells.Item(20, 13).Value = 115.0
$ws2.Cells.Item(20, 14).Value = 115.0
$ws2.Cells.Item(20, 15).Value = 115.0
$ws2.Cells.Item(20, 16).Value = 115.0
# Row 21
$ws2.Cells.Item(21, 1).Value = 45504.94097222222
$ws2.Cells.Item(21, 2).Value = 0.008
$ws2.Cells.Item(21, 3).Value = 0.004
$ws2.Cells.Item(21, 4).Value = 0.006
$ws2.Cells.Item(21, 5).Value = 0.004
$ws2.Cells.Item(21, 6).Value = 0.007
$ws2.Cells.Item(21, 7).Value = 93.61
$ws2.Cells.Item(21, 8).Value = 93.62
$ws2.Cells.Item(21, 9).Value = 93.65
$ws2.Cells.Item(21, 10).Value = 93.65
$ws2.Cells.Item(21, 11).Value = 93.65
$ws2.Cells.Item(21, 12).Value = 93.64
$ws2.Cells.Item(21, 13).Value = 93.65
$ws2.Cells.Item(21, 14).Value = 93.65
$ws2.Cells.Item(21, 15).Value = 93.63
$ws2.Cells.Item(21, 16).Value = 93.64
# Row 22
$ws2.Cells.Item(22, 1).Value = 45509.60277777778
$ws2.Cells.Item(22, 2).Value = 0.002
$ws2.Cells.Item(22, 3).Value = 0.006
$ws2.Cells.Item(22, 4).Value = 0.005
$ws2.Cells.Item(22, 5).Value = 0.004
$ws2.Cells.Item(22, 6).Value = 0.002
$ws2.Cells.Item(22, 7).Value = 62.62
$ws2.Cells.Item(22, 8).Value = 62.6
$ws2.Cells.Item(22, 9).Value = 62.61
$ws2.Cells.Item(22, 10).Value = 62.61
$ws2.Cells.Item(22, 11).Value = 62.61
$ws2.Cells.Item(22, 12).Value = 62.63
$ws2.Cells.Item(22, 13).Value = 62.63
$ws2.Cells.Item(22, 14).Value = 62.63
$ws2.Cells.Item(22, 15).Value = 62.61
$ws2.Cells.Item(22, 16).Value = 62.62
# Row 23
$ws2.Cells.Item(23, 1).Value = 45512.59722222222
$ws2.Cells.Item(23, 2).Value = 0.014
$ws2.Cells.Item(23, 3).Value = 0.005
$ws2.Cells.Item(23, 4).Value = 0.005
$ws2.Cells.Item(23, 5).Value = 0.003
$ws2.Cells.Item(23, 6).Value = 0.003
$ws2.Cells.Item(23, 7).Value = 49.23
$ws2.Cells.Item(23, 8).Value = 49.2
$ws2.Cells.Item(23, 9).Value = 49.24
$ws2.Cells.Item(23, 10).Value = 49.25
$ws2.Cells.Item(23, 11).Value = 49.24
$ws2.Cells.Item(23, 12).Value = 49.24
$ws2.Cells.Item(23, 13).Value = 49.25
$ws2.Cells.Item(23, 14).Value = 49.25
$ws2.Cells.Item(23, 15).Value = 49.24
$ws2.Cells.Item(23, 16).Value = 49.25
# Row 24
$ws2.Cells.Item(24, 1).Value = 45516.67638888889
$ws2.Cells.Item(24, 2).Value = 0.012
$ws2.Cells.Item(24, 3).Value = 0.01
$ws2.Cells.Item(24, 4).Value = 0.009
$ws2.Cells.Item(24, 5).Value = 0.01
$ws2.Cells.Item(24, 6).Value = 0.009
$ws2.Cells.Item(24, 7).Value = 34.16
$ws2.Cells.Item(24, 8).Value = 34.16
$ws2.Cells.Item(24, 9).Value = 34.17
$ws2.Cells.Item(24, 10).Value = 34.18
$ws2.Cells.Item(24, 11).Value = 34.17
$ws2.Cells.Item(24, 12).Value = 34.18
$ws2.Cells.Item(24, 13).Value = 34.19
$ws2.Cells.Item(24, 14).Value = 34.19
$ws2.Cells.Item(24, 15).Value = 34.18
$ws2.Cells.Item(24, 16).Value = 34.19
# Row 25
$ws2.Cells.Item(25, 1).Value = 45519.63888888889
$ws2.Cells.Item(25, 2).Value = 0.012
$ws2.Cells.Item(25, 3).Value = 0.011
$ws2.Cells.Item(25, 4).Value = 0.011
$ws2.Cells.Item(25, 5).Value = 0.011
$ws2.Cells.Item(25, 6).Value = 0.011
$ws2.Cells.Item(25, 7).Value = 26.64
$ws2.Cells.Item(25, 8).Value = 26.62
$ws2.Cells.Item(25, 9).Value = 26.62
$ws2.Cells.Item(25, 10).Value = 26.6
$ws2.Cells.Item(25, 11).Value = 26.59
$ws2.Cells.Item(25, 12).Value = 26.59
$ws2.Cells.Item(25, 13).Value = 26.58
$ws2.Cells.Item(25, 14).Value = 26.58
$ws2.Cells.Item(25, 15).Value = 26.58
$ws2.Cells.Item(25, 16).Value = 26.57
# Row 26
$ws2.Cells.Item(26, 1).Value = 45523.675
$ws2.Cells.Item(26, 2).Value = 0.007
$ws2.Cells.Item(26, 3).Value = 0.005
$ws2.Cells.Item(26, 4).Value = 0.005
$ws2.Cells.Item(26, 5).Value = 0.005
$ws2.Cells.Item(26, 6).Value = 0.005
$ws2.Cells.Item(26, 7).Value = 18.4
$ws2.Cells.Item(26, 8).Value = 18.39
$ws2.Cells.Item(26, 9).Value = 18.4
$ws2.Cells.Item(26, 10).Value = 18.41
$ws2.Cells.Item(26, 11).Value = 18.41
$ws2.Cells.Item(26, 12).Value = 18.41
$ws2.Cells.Item(26, 13).Value = 18.41
$ws2.Cells.Item(26, 14).Value = 18.41
$ws2.Cells.Item(26, 15).Value = 18.41
$ws2.Cells.Item(26, 16).Value = 18.41

# Row 22's "Prumer" cell picks up the banded-row border formatting that row 21
# already uses (cosmetic fix-up made while filling in the table) - copy formats
# only, so the SUM/AVERAGE formula and its freshly recalculated value are kept.
$ws2.Range("Q21").Copy()
$ws2.Range("Q22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: List1 scrolled to row 40 with B44 selected; ---
# --- Mereni aktivity (the saved active sheet) ends with A27 selected, topLeftCell default ---
$ws1.Activate()
$ws1.Range("B44").Select()
$excel.ActiveWindow.ScrollRow = 40

$ws2.Activate()
$ws2.Range("A27").Select()
